$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 0.1889250814332248
$ws.Range("C2").Value2 = 0.5830618892508144
$ws.Range("J2").Value2 = 0.003257328990228013
$ws.Range("P2").Value2 = 0.1530944625407166
$ws.Range("S2").Value2 = 0.07166123778501629
$ws.Range("C3").Value2 = 0.02173913043478261
$ws.Range("J3").Value2 = 0.0108695652173913
$ws.Range("P3").Value2 = 0.8695652173913043
$ws.Range("S3").Value2 = 0.09782608695652174
$ws.Range("P4").Value2 = 0.7555555555555555
$ws.Range("S4").Value2 = 0.2444444444444444
$ws.Range("B6").Value2 = 0.04739336492890995
$ws.Range("D6").Value2 = 0.009478672985781991
$ws.Range("F6").Value2 = 0.04265402843601896
$ws.Range("J6").Value2 = 0.1990521327014218
$ws.Range("O6").Value2 = 0.04265402843601896
$ws.Range("Q6").Value2 = 0.1800947867298578
$ws.Range("R6").Value2 = 0.1042654028436019
$ws.Range("S6").Value2 = 0.3744075829383886
$ws.Range("B7").Value2 = 0.09036144578313253
$ws.Range("D7").Value2 = 0.006024096385542169
$ws.Range("F7").Value2 = 0.05421686746987952
$ws.Range("J7").Value2 = 0.1204819277108434
$ws.Range("O7").Value2 = 0.01204819277108434
$ws.Range("Q7").Value2 = 0.1987951807228916
$ws.Range("R7").Value2 = 0.1506024096385542
$ws.Range("S7").Value2 = 0.3674698795180723
$ws.Range("B8").Value2 = 0.0962671905697446
$ws.Range("D8").Value2 = 0.02357563850687623
$ws.Range("F8").Value2 = 0.06679764243614932
$ws.Range("J8").Value2 = 0.08447937131630648
$ws.Range("O8").Value2 = 0.02946954813359529
$ws.Range("Q8").Value2 = 0.1669941060903733
$ws.Range("R8").Value2 = 0.1198428290766208
$ws.Range("S8").Value2 = 0.412573673870334
$ws.Range("B9").Value2 = 0.1129943502824859
$ws.Range("D9").Value2 = 0.02259887005649718
$ws.Range("F9").Value2 = 0.07344632768361582
$ws.Range("J9").Value2 = 0.06779661016949153
$ws.Range("O9").Value2 = 0.02824858757062147
$ws.Range("Q9").Value2 = 0.1807909604519774
$ws.Range("R9").Value2 = 0.1694915254237288
$ws.Range("S9").Value2 = 0.3446327683615819
$ws.Range("B10").Value2 = 0.1166407465007776
$ws.Range("D10").Value2 = 0.02021772939346812
$ws.Range("E10").Value2 = 0.001555209953343701
$ws.Range("F10").Value2 = 0.07387247278382582
$ws.Range("J10").Value2 = 0.08242612752721618
$ws.Range("O10").Value2 = 0.01399688958009331
$ws.Range("Q10").Value2 = 0.2068429237947123
$ws.Range("R10").Value2 = 0.1283048211508554
$ws.Range("S10").Value2 = 0.3561430793157076
$ws.Range("G11").Value2 = 0.1161825726141079
$ws.Range("J11").Value2 = 0.1037344398340249
$ws.Range("K11").Value2 = 0.1701244813278008
$ws.Range("L11").Value2 = 0.6058091286307054
$ws.Range("S11").Value2 = 0.004149377593360996
$ws.Range("G12").Value2 = 0.6887417218543046
$ws.Range("J12").Value2 = 0.2317880794701987
$ws.Range("K12").Value2 = 0.006622516556291391
$ws.Range("L12").Value2 = 0.03311258278145696
$ws.Range("S12").Value2 = 0.03973509933774835
$ws.Range("G13").Value2 = 0.7777777777777778
$ws.Range("J13").Value2 = 0.2222222222222222
$ws.Range("F15").Value2 = 0.003937007874015748
$ws.Range("H15").Value2 = 0.1692913385826772
$ws.Range("I15").Value2 = 0.06692913385826772
$ws.Range("J15").Value2 = 0.3700787401574803
$ws.Range("K15").Value2 = 0.03149606299212598
$ws.Range("M15").Value2 = 0.01574803149606299
$ws.Range("N15").Value2 = 0.003937007874015748
$ws.Range("O15").Value2 = 0.05511811023622047
$ws.Range("S15").Value2 = 0.2834645669291339
$ws.Range("F16").Value2 = 0.01724137931034483
$ws.Range("H16").Value2 = 0.1982758620689655
$ws.Range("I16").Value2 = 0.06465517241379311
$ws.Range("J16").Value2 = 0.4224137931034483
$ws.Range("K16").Value2 = 0.07327586206896551
$ws.Range("M16").Value2 = 0.01724137931034483
$ws.Range("N16").Value2 = 0.004310344827586207
$ws.Range("O16").Value2 = 0.08189655172413793
$ws.Range("S16").Value2 = 0.1206896551724138
$ws.Range("F17").Value2 = 0.01777777777777778
$ws.Range("H17").Value2 = 0.1911111111111111
$ws.Range("I17").Value2 = 0.07777777777777778
$ws.Range("J17").Value2 = 0.4533333333333333
$ws.Range("K17").Value2 = 0.09111111111111111
$ws.Range("M17").Value2 = 0.02
$ws.Range("N17").Value2 = 0.004444444444444444
$ws.Range("O17").Value2 = 0.04888888888888889
$ws.Range("S17").Value2 = 0.09555555555555556
$ws.Range("F18").Value2 = 0.006734006734006734
$ws.Range("H18").Value2 = 0.1986531986531987
$ws.Range("I18").Value2 = 0.07744107744107744
$ws.Range("J18").Value2 = 0.4343434343434344
$ws.Range("K18").Value2 = 0.08080808080808081
$ws.Range("M18").Value2 = 0.0303030303030303
$ws.Range("O18").Value2 = 0.08417508417508418
$ws.Range("S18").Value2 = 0.08754208754208755
$ws.Range("F19").Value2 = 0.008298755186721992
$ws.Range("H19").Value2 = 0.2257261410788382
$ws.Range("I19").Value2 = 0.07302904564315353
$ws.Range("J19").Value2 = 0.4024896265560166
$ws.Range("K19").Value2 = 0.08796680497925312
$ws.Range("M19").Value2 = 0.01659751037344398
$ws.Range("N19").Value2 = 0.0008298755186721991
$ws.Range("O19").Value2 = 0.08132780082987552
$ws.Range("S19").Value2 = 0.1037344398340249
